$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Slit1","Robo1","ECs",3,1,0.340305,1.020915,0.2413792532744959,0.2522899505114672,3,1,0.1438866666666667,0.43166,0.005528090411415381,0.005816871971504955,0.04896535210000001,0.4406881689000001,0.001334366335541346,0.001467538341822526),
    @("FAPs","Slit1","Robo1","FAPs",3,1,0.340305,1.020915,0.2413792532744959,0.2522899505114672,3,1,21.99231,65.97693,0.8449391514331273,0.889077873518381,7.484093054550001,67.35683749095,0.2039507814353145,0.2243054127107928),
    @("FAPs","Slit1","Robo1","M1",3,1,0.340305,1.020915,0.2413792532744959,0.2522899505114672,1,0.3333333333333333,0.015519,0.046557,0.0005962361703291152,0.0006273829133516105,0.005281193295,0.047530739655,0.000143919041569287,0.0001582824041612179),
    @("FAPs","Slit1","Robo1","sCs",3,1,0.340305,1.020915,0.2413792532744959,0.2522899505114672,2,1,3.876561,7.753121999999999,0.1489365219851282,0.1044778715967623,1.319213091105,7.91527854663,0.03595018646207079,0.02635871705469059),
    @("M1","Slit1","Robo1","ECs",3,1,0.472982,1.418946,0.3354874068035369,0.3506519309819567,3,1,0.1438866666666667,0.43166,0.005528090411415381,0.005816871971504955,0.06805580337333335,0.6125022303600001,0.001854604716701244,0.002039697389083034),
    @("M1","Slit1","Robo1","FAPs",3,1,0.472982,1.418946,0.3354874068035369,0.3506519309819567,3,1,21.99231,65.97693,0.8449391514331273,0.889077873518381,10.40196676842,93.61770091577999,0.2834664448210808,0.3117568731425521),
    @("M1","Slit1","Robo1","M1",3,1,0.472982,1.418946,0.3354874068035369,0.3506519309819567,1,0.3333333333333333,0.015519,0.046557,0.0005962361703291152,0.0006273829133516105,0.007340207658,0.066061868922,0.0002000297266261868,0.0002199930300318278),
    @("M1","Slit1","Robo1","sCs",3,1,0.472982,1.418946,0.3354874068035369,0.3506519309819567,2,1,3.876561,7.753121999999999,0.1489365219851282,0.1044778715967623,1.833543574902,11.001261449412,0.04996632753912863,0.03663536742028964),
    @("M2","Slit1","Robo1","ECs",3,1,0.4136363333333333,1.240909,0.2933933655608953,0.3066551771687498,3,1,0.1438866666666667,0.43166,0.005528090411415381,0.005816871971504955,0.05951675321555557,0.53565077894,0.001621905050930073,0.001783773904989787),
    @("M2","Slit1","Robo1","FAPs",3,1,0.4136363333333333,1.240909,0.2933933655608953,0.3066551771687498,3,1,21.99231,65.97693,0.8449391514331273,0.889077873518381,9.09681846993,81.87136622937,0.2478995413331322,0.2726403328205945),
    @("M2","Slit1","Robo1","M1",3,1,0.4136363333333333,1.240909,0.2933933655608953,0.3066551771687498,1,0.3333333333333333,0.015519,0.046557,0.0005962361703291152,0.0006273829133516105,0.006419222257,0.05777300031300001,0.0001749317366819983,0.0001923902184464846),
    @("M2","Slit1","Robo1","sCs",3,1,0.4136363333333333,1.240909,0.2933933655608953,0.3066551771687498,2,1,3.876561,7.753121999999999,0.1489365219851282,0.1044778715967623,1.603486477983,9.620918867898,0.04369698744015104,0.03203868022471905),
    @("sCs","Slit1","Robo1","ECs",2,1,0.182912,0.365824,0.129739974361072,0.09040294133782634,3,1,0.1438866666666667,0.43166,0.005528090411415381,0.005816871971504955,0.02631859797333334,0.15791158784,0.0007172143082427193,0.0005258623356096087),
    @("sCs","Slit1","Robo1","FAPs",2,1,0.182912,0.365824,0.129739974361072,0.09040294133782634,3,1,21.99231,65.97693,0.8449391514331273,0.889077873518381,4.02265740672,24.13594444032,0.1096223838435998,0.08037525484444158),
    @("sCs","Slit1","Robo1","M1",2,1,0.182912,0.365824,0.129739974361072,0.09040294133782634,1,0.3333333333333333,0.015519,0.046557,0.0005962361703291152,0.0006273829133516105,0.002838611328,0.017031667968,[double]"7.735566545164314E-05",[double]"5.671726071208023E-05"),
    @("sCs","Slit1","Robo1","sCs",2,1,0.182912,0.365824,0.129739974361072,0.09040294133782634,2,1,3.876561,7.753121999999999,0.1489365219851282,0.1044778715967623,0.709069525632,2.836278102528,0.01932302054377777,0.009445106897063058)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
